$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header row: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404" ---
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"
$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# --- Turn the used range into an Excel Table ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U77"), 0, 1)
$tbl.Name = "Table1"

# --- Freeze the header row ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
